$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D:E columns for data rows so numeric-looking strings
# (e.g. "417.45") are stored as text, matching the source data which uses
# inline/shared strings rather than numbers.
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = "62.894.24"
$ws.Range("E2").Value = "  +6.25%  "

$ws.Range("D3").Value = "3.482.20"
$ws.Range("E3").Value = "  +5.04%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "417.45"
$ws.Range("E5").Value = "  +3.28%  "

$ws.Range("D6").Value = "132.75"
$ws.Range("E6").Value = "  +20.77%  "

$ws.Range("D7").Value = "3.473.45"
$ws.Range("E7").Value = "  +4.87%  "

$ws.Range("D8").Value = "0.596"
$ws.Range("E8").Value = "  +1.95%  "

$ws.Range("D9").Value = "1.00"
$ws.Range("E9").Value = "  +0.10%  "

$ws.Range("D10").Value = "0.697"
$ws.Range("E10").Value = "  +9.97%  "

$ws.Range("D11").Value = "0.127"
$ws.Range("E11").Value = "  +29.31%  "

$ws.Range("D12").Value = "44.43"
$ws.Range("E12").Value = "  +12.19%  "

$ws.Range("E13").Value = "  +0.47%  "

$ws.Range("D14").Value = "4.028.50"
$ws.Range("E14").Value = "  +4.57%  "

$ws.Range("D15").Value = "8.88"
$ws.Range("E15").Value = "  +6.01%  "

$ws.Range("D16").Value = "20.37"
$ws.Range("E16").Value = "  +5.73%  "

$ws.Range("D17").Value = "3.479.04"
$ws.Range("E17").Value = "  +4.89%  "

$ws.Range("D18").Value = "62.818.68"
$ws.Range("E18").Value = "  +6.43%  "

$ws.Range("D19").Value = "1.05"
$ws.Range("E19").Value = "  +1.17%  "

$ws.Range("D20").Value = "11.02"
$ws.Range("E20").Value = "  +3.02%  "

$ws.Range("D21").Value = "0.0000137"
$ws.Range("E21").Value = "  +25.27%  "

$ws.Range("D22").Value = "3.39"
$ws.Range("E22").Value = "  +3.03%  "

$ws.Range("E23").Value = "  +3.04%  "

$ws.Range("D24").Value = "82.51"
$ws.Range("E24").Value = "  +10.22%  "

$ws.Range("D25").Value = "317.02"
$ws.Range("E25").Value = "  +4.22%  "

$ws.Range("D26").Value = "3.23"
$ws.Range("E26").Value = "  +1.28%  "

$ws.Range("D27").Value = "30.91"
$ws.Range("E27").Value = "  +8.83%  "

$ws.Range("D28").Value = "8.17"
$ws.Range("E28").Value = "  +4.23%  "

$ws.Range("D29").Value = "7.84"
$ws.Range("E29").Value = "  +7.92%  "

$ws.Range("D30").Value = "0.179"
$ws.Range("E30").Value = "  +5.15%  "

$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "0.121"
$ws.Range("E31").Value = "  +8.55%  "

$ws.Range("B32").Value = "LEO"
$ws.Range("C32").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D32").Value = "4.40"
$ws.Range("E32").Value = "  -1.16%  "

$ws.Range("D33").Value = "44.76"
$ws.Range("E33").Value = "  +11.97%  "

$ws.Range("B34").Value = "Cosmos"
$ws.Range("C34").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D34").Value = "11.93"
$ws.Range("E34").Value = "  +5.16%  "

$ws.Range("B35").Value = "Toncoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D35").Value = "2.66"
$ws.Range("E35").Value = "  +24.94%  "

$ws.Range("E36").Value = "  +0.09%  "

$ws.Range("D37").Value = "0.0498"
$ws.Range("E37").Value = "  -4.92%  "

$ws.Range("D38").Value = "52.56"
$ws.Range("E38").Value = "  +1.49%  "

$ws.Range("D39").Value = "3.61"
$ws.Range("E39").Value = "  +4.33%  "

$ws.Range("D40").Value = "0.996"
$ws.Range("E40").Value = "  -0.54%  "

$ws.Range("D41").Value = "3.05"
$ws.Range("E41").Value = "  -6.42%  "

$ws.Range("D42").Value = "2.03"
$ws.Range("E42").Value = "  +7.90%  "

$ws.Range("E43").Value = "  +3.55%  "

$ws.Range("D44").Value = "137.74"
$ws.Range("E44").Value = "  +0.25%  "

$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "4.08"
$ws.Range("E45").Value = "  +4.85%  "

$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D46").Value = "17.50"
$ws.Range("E46").Value = "  +5.83%  "

$ws.Range("E47").Value = "  +4.60%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "22.78"
$ws.Range("E48").Value = "  +2.59%  "

$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "2.26"
$ws.Range("E49").Value = "  -2.45%  "

$ws.Range("D50").Value = "2.249.58"
$ws.Range("E50").Value = "  +3.78%  "

$ws.Range("D51").Value = "3.822.07"
$ws.Range("E51").Value = "  +4.70%  "

# Restore default cell style (the text number format was only needed to
# prevent Excel from auto-converting numeric-looking text into numbers).
$priceVolRange.Style = "Normal"
